# Scheduled market-data refresh for Mandragora_Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# on each job sheet with the latest pulled market values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow
$ws.Range("H12").Value = 117.5
$ws.Range("I12").Value = 110
$ws.Range("J12").Value = 125
$ws.Range("K12").Value = 110
$ws.Range("L12").Value = 125
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = -465

# Row 19: Unbreak My Heart
$ws.Range("H19").Value = 610.5714
$ws.Range("I19").Value = 454.30768
$ws.Range("J19").Value = 746
$ws.Range("K19").Value = 454.30768
$ws.Range("L19").Value = 746
$ws.Range("M19").Value = -279.30768
$ws.Range("N19").Value = -1096

# Row 51: A Bile Business
$ws.Range("H51").Value = 9096759
$ws.Range("I51").Value = 18187218
$ws.Range("J51").Value = 6300.4
$ws.Range("K51").Value = 18187218
$ws.Range("L51").Value = 6300.4
$ws.Range("M51").Value = -18186734
$ws.Range("N51").Value = -7268.4

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 603.5
$ws.Range("I107").Value = 616.2222
$ws.Range("J107").Value = 565.3333
$ws.Range("K107").Value = 616.2222
$ws.Range("L107").Value = 565.3333
$ws.Range("M107").Value = 1303.7778
$ws.Range("N107").Value = -4405.3333

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 4992.125
$ws.Range("I132").Value = 3835.8684
$ws.Range("J132").Value = 7433.1113
$ws.Range("K132").Value = 11507.6052
$ws.Range("L132").Value = 22299.3339
$ws.Range("M132").Value = -8977.6052
$ws.Range("N132").Value = -27359.3339

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 7161.6353
$ws.Range("I32").Value = 7803.0547
$ws.Range("J32").Value = 5304.8945
$ws.Range("K32").Value = 7803.0547
$ws.Range("L32").Value = 5304.8945
$ws.Range("M32").Value = -7516.0547
$ws.Range("N32").Value = -5878.8945

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4524.159
$ws.Range("I132").Value = 2102.318
$ws.Range("J132").Value = 6946
$ws.Range("K132").Value = 6306.954000000001
$ws.Range("L132").Value = 20838
$ws.Range("M132").Value = -3776.954000000001
$ws.Range("N132").Value = -25898

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 1214.6154
$ws.Range("I22").Value = 1354.3334
$ws.Range("J22").Value = 900.25
$ws.Range("K22").Value = 1354.3334
$ws.Range("L22").Value = 900.25
$ws.Range("M22").Value = -1181.3334
$ws.Range("N22").Value = -1246.25

# Row 43: Don't Fear the Reaper
$ws.Range("H43").Value = 116000
$ws.Range("J43").Value = 116000
$ws.Range("L43").Value = 116000
$ws.Range("N43").Value = -116362

# Row 64: With Bearings Straight
$ws.Range("H64").Value = 1072.2858
$ws.Range("I64").Value = 1006
$ws.Range("K64").Value = 1006
$ws.Range("M64").Value = -781

# Row 67: Bearing the Brunt (L)
$ws.Range("H67").Value = 1072.2858
$ws.Range("I67").Value = 1006
$ws.Range("K67").Value = 1006
$ws.Range("M67").Value = -226

# Row 75: I Saw the Pine
$ws.Range("H75").Value = 5000
$ws.Range("I75").Value = 5000
$ws.Range("K75").Value = 5000
$ws.Range("M75").Value = -4064

# Row 78: I Came, I Sawed, I Conquered (L)
$ws.Range("H78").Value = 5000
$ws.Range("I78").Value = 5000
$ws.Range("K78").Value = 15000
$ws.Range("M78").Value = -10320

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent
$ws.Range("H7").Value = 4650.609
$ws.Range("I7").Value = 9120.637000000001
$ws.Range("K7").Value = 9120.637000000001
$ws.Range("M7").Value = -9007.637000000001

# Row 31: Wall Not Found
$ws.Range("H31").Value = 2386.0952
$ws.Range("I31").Value = 1750.25
$ws.Range("J31").Value = 4420.8
$ws.Range("K31").Value = 1750.25
$ws.Range("L31").Value = 4420.8
$ws.Range("M31").Value = -1455.25
$ws.Range("N31").Value = -5010.8

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2386.0952
$ws.Range("I34").Value = 1750.25
$ws.Range("J34").Value = 4420.8
$ws.Range("K34").Value = 1750.25
$ws.Range("L34").Value = 4420.8
$ws.Range("M34").Value = -1548.25
$ws.Range("N34").Value = -4824.8

# Row 55: Ready for a Rematch
$ws.Range("H55").Value = 8633
$ws.Range("I55").Value = 8633
$ws.Range("K55").Value = 8633
$ws.Range("M55").Value = -8318

# Row 94: Beech, Please
$ws.Range("H94").Value = 984.625
$ws.Range("J94").Value = 984.625
$ws.Range("L94").Value = 984.625
$ws.Range("N94").Value = -1886.625

# Row 99: O Pine
$ws.Range("H99").Value = 4191.778
$ws.Range("I99").Value = 2422.4
$ws.Range("J99").Value = 6403.5
$ws.Range("K99").Value = 2422.4
$ws.Range("L99").Value = 6403.5
$ws.Range("M99").Value = -924.4000000000001
$ws.Range("N99").Value = -9399.5

# Row 107: Built to Last
$ws.Range("H107").Value = 514.5454999999999
$ws.Range("I107").Value = 231.27272
$ws.Range("J107").Value = 797.8182
$ws.Range("K107").Value = 231.27272
$ws.Range("L107").Value = 797.8182
$ws.Range("M107").Value = 1688.72728
$ws.Range("N107").Value = -4637.8182

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4191.778
$ws.Range("I126").Value = 2422.4
$ws.Range("J126").Value = 6403.5
$ws.Range("K126").Value = 7267.200000000001
$ws.Range("L126").Value = 19210.5
$ws.Range("M126").Value = -4797.200000000001
$ws.Range("N126").Value = -24150.5

# Row 138: Bow Out
$ws.Range("H138").Value = 45540
$ws.Range("J138").Value = 45540
$ws.Range("L138").Value = 45540
$ws.Range("N138").Value = -55820

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Range("H34").Value = 20112.092
$ws.Range("J34").Value = 22073.3
$ws.Range("L34").Value = 66219.89999999999
$ws.Range("N34").Value = -66387.89999999999

# Row 68: Such a Butter Face
$ws.Range("H68").Value = 627.2222
$ws.Range("I68").Value = 521
$ws.Range("J68").Value = 999
$ws.Range("K68").Value = 1563
$ws.Range("L68").Value = 2997
$ws.Range("M68").Value = -752
$ws.Range("N68").Value = -4619

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 627.2222
$ws.Range("I71").Value = 521
$ws.Range("J71").Value = 999
$ws.Range("K71").Value = 4689
$ws.Range("L71").Value = 8991
$ws.Range("M71").Value = -633
$ws.Range("N71").Value = -17103

# Row 107: Slippery Service
$ws.Range("H107").Value = 18519092
$ws.Range("I107").Value = 125000200
$ws.Range("J107").Value = 638.7826
$ws.Range("K107").Value = 375000600
$ws.Range("L107").Value = 1916.3478
$ws.Range("M107").Value = -374998680
$ws.Range("N107").Value = -5756.3478

# Row 132: More Mezcal
$ws.Range("H132").Value = 3861.5557
$ws.Range("I132").Value = 1750
$ws.Range("K132").Value = 15750
$ws.Range("M132").Value = -13220

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 3057
$ws.Range("I7").Value = 2825
$ws.Range("J7").Value = 3985
$ws.Range("K7").Value = 2825
$ws.Range("L7").Value = 3985
$ws.Range("M7").Value = -2713
$ws.Range("N7").Value = -4209

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 2473.9565
$ws.Range("I93").Value = 2213.1667
$ws.Range("J93").Value = 3412.8
$ws.Range("K93").Value = 2213.1667
$ws.Range("L93").Value = 3412.8
$ws.Range("M93").Value = -965.1667000000002
$ws.Range("N93").Value = -5908.8

# Row 126: Battered Books
$ws.Range("H126").Value = 3057
$ws.Range("I126").Value = 2825
$ws.Range("J126").Value = 3985
$ws.Range("K126").Value = 8475
$ws.Range("L126").Value = 11955
$ws.Range("M126").Value = -6005
$ws.Range("N126").Value = -16895

$ws = $wb.Worksheets.Item("WVR")
# Row 137: Traditional Trousers
$ws.Range("H137").Value = 32887.855
$ws.Range("J137").Value = 32887.855
$ws.Range("L137").Value = 32887.855
$ws.Range("N137").Value = -43087.855
